# Update currentAveragePrice / LevePrice / LeveProfit columns (H-N) for items
# whose market-board snapshot moved since the last scheduled refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 4226.0435
$ws.Range("I116").Value = 2224.9167
$ws.Range("J116").Value = 6409.091
$ws.Range("K116").Value = 2224.9167
$ws.Range("L116").Value = 6409.091
$ws.Range("M116").Value = 1217.0833
$ws.Range("N116").Value = -13293.091

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 525.6429000000001
$ws.Range("I125").Value = 290.25
$ws.Range("J125").Value = 619.8
$ws.Range("K125").Value = 2612.25
$ws.Range("L125").Value = 5578.2
$ws.Range("M125").Value = -152.25
$ws.Range("N125").Value = -10498.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1879.919
$ws.Range("J129").Value = 1909.862
$ws.Range("L129").Value = 5729.586
$ws.Range("N129").Value = -15729.586

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1620.8889
$ws.Range("I2").Value = 1614.56
$ws.Range("J2").Value = 1700
$ws.Range("K2").Value = 1614.56
$ws.Range("L2").Value = 1700
$ws.Range("M2").Value = -1501.56
$ws.Range("N2").Value = -1926

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9942.478999999999
$ws.Range("I32").Value = 9154.282999999999
$ws.Range("J32").Value = 15197.111
$ws.Range("K32").Value = 9154.282999999999
$ws.Range("L32").Value = 15197.111
$ws.Range("M32").Value = -8867.282999999999
$ws.Range("N32").Value = -15771.111

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2343.353
$ws.Range("I45").Value = 2162.4
$ws.Range("J45").Value = 2601.8572
$ws.Range("K45").Value = 2162.4
$ws.Range("L45").Value = 2601.8572
$ws.Range("M45").Value = -1785.4
$ws.Range("N45").Value = -3355.8572

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1793.55
$ws.Range("I74").Value = 1744.2069
$ws.Range("J74").Value = 1923.6364
$ws.Range("K74").Value = 1744.2069
$ws.Range("L74").Value = 1923.6364
$ws.Range("M74").Value = -870.2068999999999
$ws.Range("N74").Value = -3671.6364

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1793.55
$ws.Range("I77").Value = 1744.2069
$ws.Range("J77").Value = 1923.6364
$ws.Range("K77").Value = 8721.0345
$ws.Range("L77").Value = 9618.182000000001
$ws.Range("M77").Value = -4353.0345
$ws.Range("N77").Value = -18354.182

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1705.3914
$ws.Range("I110").Value = 1742.6842
$ws.Range("J110").Value = 1528.25
$ws.Range("K110").Value = 1742.6842
$ws.Range("L110").Value = 1528.25
$ws.Range("M110").Value = 302.3158000000001
$ws.Range("N110").Value = -5618.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1620.8889
$ws.Range("I116").Value = 1614.56
$ws.Range("J116").Value = 1700
$ws.Range("K116").Value = 1614.56
$ws.Range("L116").Value = 1700
$ws.Range("M116").Value = 679.4400000000001
$ws.Range("N116").Value = -6288

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 11365307
$ws.Range("I132").Value = 14706939
$ws.Range("J132").Value = 3757.4
$ws.Range("K132").Value = 44120817
$ws.Range("L132").Value = 11272.2
$ws.Range("M132").Value = -44118287
$ws.Range("N132").Value = -16332.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 20000
$ws.Range("J2").Value = 20000
$ws.Range("L2").Value = 20000
$ws.Range("N2").Value = -20226

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1620.8889
$ws.Range("I3").Value = 1614.56
$ws.Range("J3").Value = 1700
$ws.Range("K3").Value = 1614.56
$ws.Range("L3").Value = 1700
$ws.Range("M3").Value = -1500.56
$ws.Range("N3").Value = -1928

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 18400
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 18400
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 18400
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value = -23892

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1574.196
$ws.Range("I107").Value = 1517.8485
$ws.Range("J107").Value = 1677.5
$ws.Range("K107").Value = 1517.8485
$ws.Range("L107").Value = 1677.5
$ws.Range("M107").Value = 402.1514999999999
$ws.Range("N107").Value = -5517.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2339.0708
$ws.Range("I134").Value = 962.1539
$ws.Range("J134").Value = 3862.468
$ws.Range("K134").Value = 2886.4617
$ws.Range("L134").Value = 11587.404
$ws.Range("M134").Value = -351.4616999999998
$ws.Range("N134").Value = -16657.404

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1327.4807
$ws.Range("I58").Value = 1035.1282
$ws.Range("J58").Value = 2204.5386
$ws.Range("K58").Value = 1035.1282
$ws.Range("L58").Value = 2204.5386
$ws.Range("M58").Value = -832.1282000000001
$ws.Range("N58").Value = -2610.5386

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 24391.115
$ws.Range("I132").Value = 1263.9387
$ws.Range("J132").Value = 118827.086
$ws.Range("K132").Value = 3791.8161
$ws.Range("L132").Value = 356481.258
$ws.Range("M132").Value = -1261.8161
$ws.Range("N132").Value = -361541.258

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1110.2727
$ws.Range("I134").Value = 825.9286
$ws.Range("J134").Value = 2702.6
$ws.Range("K134").Value = 2477.7858
$ws.Range("L134").Value = 8107.799999999999
$ws.Range("M134").Value = 57.21420000000035
$ws.Range("N134").Value = -13177.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1327.4807
$ws.Range("I136").Value = 1035.1282
$ws.Range("J136").Value = 2204.5386
$ws.Range("K136").Value = 3105.3846
$ws.Range("L136").Value = 6613.6158
$ws.Range("M136").Value = -555.3846000000003
$ws.Range("N136").Value = -11713.6158

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 11016.211
$ws.Range("I107").Value = 11548.444
$ws.Range("J107").Value = 10537.2
$ws.Range("K107").Value = 34645.33199999999
$ws.Range("L107").Value = 31611.6
$ws.Range("M107").Value = -32725.33199999999
$ws.Range("N107").Value = -35451.60000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 753600
$ws.Range("I120").Value = 1000800
$ws.Range("J120").Value = 12000
$ws.Range("K120").Value = 3002400
$ws.Range("L120").Value = 36000
$ws.Range("M120").Value = -2997562
$ws.Range("N120").Value = -45676

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1686.5
$ws.Range("I102").Value = 1824.8
$ws.Range("J102").Value = 995
$ws.Range("K102").Value = 1824.8
$ws.Range("L102").Value = 995
$ws.Range("M102").Value = -202.8
$ws.Range("N102").Value = -4239

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 14414.429
$ws.Range("J123").Value = 14414.429
$ws.Range("L123").Value = 14414.429
$ws.Range("N123").Value = -19314.429

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H134").Value = 25353.25
$ws.Range("J134").Value = 25353.25
$ws.Range("L134").Value = 76059.75
$ws.Range("N134").Value = -81129.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 49881.5
$ws.Range("J136").Value = 49881.5
$ws.Range("L136").Value = 149644.5
$ws.Range("N136").Value = -154744.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1305.2444
$ws.Range("I136").Value = 1042.3055
$ws.Range("J136").Value = 2357
$ws.Range("K136").Value = 3126.9165
$ws.Range("L136").Value = 7071
$ws.Range("M136").Value = -576.9164999999998
$ws.Range("N136").Value = -12171

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 4235999.5
$ws.Range("I5").Value = 5500000
$ws.Range("J5").Value = 4067466.2
$ws.Range("K5").Value = 5500000
$ws.Range("L5").Value = 4067466.2
$ws.Range("M5").Value = -5499888
$ws.Range("N5").Value = -4067690.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 8334033.5
$ws.Range("I107").Value = 540.4
$ws.Range("J107").Value = 14286529
$ws.Range("K107").Value = 1621.2
$ws.Range("L107").Value = 42859587
$ws.Range("M107").Value = 298.8000000000002
$ws.Range("N107").Value = -42863427

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H118").Value = 43384
$ws.Range("J118").Value = 43384
$ws.Range("L118").Value = 43384
$ws.Range("N118").Value = -46698

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 981908.0600000001
$ws.Range("I126").Value = 1549358.9
$ws.Range("J126").Value = 1765.6364
$ws.Range("K126").Value = 4648076.699999999
$ws.Range("L126").Value = 5296.9092
$ws.Range("M126").Value = -4645606.699999999
$ws.Range("N126").Value = -10236.9092

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1231.4606
$ws.Range("I132").Value = 1007.35596
$ws.Range("J132").Value = 2009.2354
$ws.Range("K132").Value = 3022.06788
$ws.Range("L132").Value = 6027.706200000001
$ws.Range("M132").Value = -492.0678800000001
$ws.Range("N132").Value = -11087.7062
